# Update crypto price/volume figures (scheduled data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds free-form text (e.g. "1.781.01", "0.0₃0707"),
# not numeric values, so force text formatting before writing the new
# figures to prevent Excel's automatic number coercion from mangling them
# (dropping trailing zeros, switching "." thousands separators to floats, etc).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.958.99"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.558.55"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "22.12"
$ws.Range("E8").Value = "  +0.22%  "

$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").Value = "1.780.16"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "1.550.34"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").Value = "61.90"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "26.950.06"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +1.63%  "

$ws.Range("D19").Value = "215.92"

$ws.Range("D20").Value = "7.38"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  +1.43%  "

$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("E24").Value = "  -1.29%  "

$ws.Range("D25").Value = "152.83"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("E27").Value = "  +0.97%  "

$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("E33").Value = "  +2.72%  "

$ws.Range("D34").Value = "1.424.33"
$ws.Range("E34").Value = "  +0.14%  "

$ws.Range("E35").Value = "  +0.75%  "

$ws.Range("E36").Value = "  +8.16%  "

$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("E39").Value = "  +2.17%  "

$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("E43").Value = "  +2.04%  "

$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "64.54"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "1.693.70"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").Value = "87.18"
$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  +4.21%  "

$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").Value = "0.0960"
$ws.Range("E51").Value = "  -0.07%  "
